$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "42.881.64"
$ws.Range("E2").Value = "  +0.00%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.278.37"
$ws.Range("E3").Value = "  -0.15%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.12%  "

# Row 5: BNB
$s = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.54"
$ws.Range("D5").Style = $s
$ws.Range("E5").Value = "  -1.27%  "

# Row 6: XRP
$s = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.632"
$ws.Range("D6").Style = $s
$ws.Range("E6").Value = "  +0.33%  "

# Row 7: Solana
$s = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "78.80"
$ws.Range("D7").Style = $s
$ws.Range("E7").Value = "  +8.09%  "

# Row 8: USDC
$ws.Range("E8").Value = "  -0.03%  "

# Row 9: Cardano
$s = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.646"
$ws.Range("D9").Style = $s
$ws.Range("E9").Value = "  -1.29%  "

# Row 10: Avalanche
$s = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.46"
$ws.Range("D10").Style = $s
$ws.Range("E10").Value = "  +6.85%  "

# Row 11: Dogecoin
$s = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0972"
$ws.Range("D11").Style = $s
$ws.Range("E11").Value = "  -0.51%  "

# Row 12: Polkadot
$s = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.34"
$ws.Range("D12").Style = $s
$ws.Range("E12").Value = "  -0.44%  "

# Row 13: TRON
$s = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.105"
$ws.Range("D13").Style = $s
$ws.Range("E13").Value = "  +0.50%  "

# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.618.45"
$ws.Range("E14").Value = "  -0.27%  "

# Row 15: Chainlink
$s = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.07"
$ws.Range("D15").Style = $s
$ws.Range("E15").Value = "  +0.36%  "

# Row 16: Polygon
$s = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.866"
$ws.Range("D16").Style = $s
$ws.Range("E16").Value = "  -2.54%  "

# Row 17: WrappedEther
$ws.Range("D17").Value = "2.284.01"
$ws.Range("E17").Value = "  +0.02%  "

# Row 18: WrappedBTC
$ws.Range("D18").Value = "42.772.95"
$ws.Range("E18").Value = "  -0.15%  "

# Row 19: ShibaInu
$ws.Range("D19").Value = "0.0₃0994"
$ws.Range("E19").Value = "  -2.23%  "

# Row 20: Uniswap
$s = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.23"
$ws.Range("D20").Style = $s
$ws.Range("E20").Value = "  -1.87%  "

# Row 21: Litecoin
$s = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.09"
$ws.Range("D21").Style = $s
$ws.Range("E21").Value = "  -1.89%  "

# Row 22: BitcoinCash
$s = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "233.53"
$ws.Range("D22").Style = $s
$ws.Range("E22").Value = "  -1.56%  "

# Row 23: ImmutableX
$s = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.16"
$ws.Range("D23").Style = $s
$ws.Range("E23").Value = "  +0.86%  "

# Row 24: WEMIXToken
$ws.Range("E24").Value = "  -2.63%  "

# Row 25: Dai
$ws.Range("E25").Value = "  +0.02%  "

# Row 26: Cosmos
$s = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.28"
$ws.Range("D26").Style = $s
$ws.Range("E26").Value = "  -3.50%  "

# Row 27: PancakeSwap
$s = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.34"
$ws.Range("D27").Style = $s
$ws.Range("E27").Value = "  -4.74%  "

# Row 28: Toncoin
$ws.Range("E28").Value = "  +2.14%  "

# Row 29: Monero
$s = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "168.24"
$ws.Range("D29").Style = $s
$ws.Range("E29").Value = "  +0.11%  "

# Row 30: EthereumClassic
$s = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.86"
$ws.Range("D30").Style = $s
$ws.Range("E30").Value = "  -0.97%  "

# Row 31: InternetComputer(DFINITY)
$s = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.51"
$ws.Range("D31").Style = $s
$ws.Range("E31").Value = "  +2.29%  "

# Row 32: Hedera
$s = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0860"
$ws.Range("D32").Style = $s
$ws.Range("E32").Value = "  +5.67%  "

# Row 33: Kaspa
$ws.Range("E33").Value = "  -4.70%  "

# Row 34: InjectiveProtocol
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$s = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.42"
$ws.Range("D34").Style = $s
$ws.Range("E34").Value = "  -1.32%  "

# Row 35: Stellar
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$s = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.127"
$ws.Range("D35").Style = $s
$ws.Range("E35").Value = "  +0.20%  "

# Row 36: RenderToken
$ws.Range("E36").Value = "  -4.22%  "

# Row 37: Filecoin
$ws.Range("E37").Value = "  -0.56%  "

# Row 38: VeChain
$s = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0303"
$ws.Range("D38").Style = $s
$ws.Range("E38").Value = "  -2.15%  "

# Row 39: Celestia
$s = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.69"
$ws.Range("D39").Style = $s
$ws.Range("E39").Value = "  +2.92%  "

# Row 40: LidoDAOToken
$ws.Range("E40").Value = "  -3.00%  "

# Row 41: THORChain
$s = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.87"
$ws.Range("D41").Style = $s
$ws.Range("E41").Value = "  -1.83%  "

# Row 42: Aave
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$s = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "112.75"
$ws.Range("D42").Style = $s
$ws.Range("E42").Value = "  +16.49%  "

# Row 43: Algorand
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$s = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.208"
$ws.Range("D43").Style = $s
$ws.Range("E43").Value = "  -1.65%  "

# Row 44: MultiversX
$s = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "61.17"
$ws.Range("D44").Style = $s
$ws.Range("E44").Value = "  -0.36%  "

# Row 45: FraxShare
$ws.Range("E45").Value = "  -3.01%  "

# Row 46: Cronos
$s = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.101"
$ws.Range("D46").Style = $s
$ws.Range("E46").Value = "  -1.94%  "

# Row 47: BinanceUSD
$ws.Range("B47").Value = "BinanceUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$s = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").Style = $s
$ws.Range("E47").Value = "  -0.09%  "

# Row 48: FTXToken
$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$s = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.54"
$ws.Range("D48").Style = $s
$ws.Range("E48").Value = "  -8.61%  "

# Row 49: ARBITRUM
$ws.Range("E49").Value = "  -2.97%  "

# Row 50: TrustWalletToken
$ws.Range("E50").Value = "  -3.14%  "

# Row 51: SynthetixNetwork
$s = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.24"
$ws.Range("D51").Style = $s
$ws.Range("E51").Value = "  -1.01%  "
